$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a 2D act camera offset/rotation for the village scene (row 2)
$ws.Range("J2").Value = "0,4.2,5.5"
$ws.Range("K2").Value = "25,180"

[void]$ws.Range("K2").Select()
